$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original text formatting
# (values like "113.52" or "2.40" would otherwise be auto-converted to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.803.59"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.249.43"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "113.52"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "295.05"
$ws.Range("E6").Value = "  +6.66%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "44.12"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").Value = "0.0924"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "54.49"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "8.90"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "1.06"
$ws.Range("E14").Value = "  +22.21%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "15.10"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "2.580.70"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "2.242.10"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "42.678.38"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  +7.17%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("D22").Value = "74.88"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("D23").Value = "3.36"
$ws.Range("E23").Value = "  +10.98%  "
$ws.Range("D26").Value = "8.99"
$ws.Range("E26").Value = "  -2.92%  "
$ws.Range("D28").Value = "11.54"
$ws.Range("E28").Value = "  -5.39%  "
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").Value = "175.42"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("D31").Value = "37.67"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("D32").Value = "21.89"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "0.0886"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "5.71"
$ws.Range("E35").Value = "  +2.50%  "
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +9.02%  "
$ws.Range("D37").Value = "4.27"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "0.0376"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("D41").Value = "2.42"
$ws.Range("E41").Value = "  -5.40%  "
$ws.Range("D42").Value = "71.88"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "0.231"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "12.49"
$ws.Range("E45").Value = "  -5.37%  "
$ws.Range("D46").Value = "1.32"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "5.49"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").Value = "1.30"
$ws.Range("E48").Value = "  +2.41%  "
$ws.Range("D49").Value = "105.20"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("D50").Value = "8.60"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("D51").Value = "0.0981"
$ws.Range("E51").Value = "  -0.57%  "

# Row 24/25 swap: BitcoinCash <-> ImmutableX with updated values
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").Value = "250.74"
$ws.Range("E25").Value = "  +8.33%  "
